$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 24.09000000000033
$ws.Range("G2").Value = [double]"7.488454301096681e-13"
$ws.Range("H2").Value = [double]"3.34989767883805e-12"
$ws.Range("K2").Value = 44.39812637908831
$ws.Range("L2").Value = "[29.938956254766055, 58.857296503410566]"
$ws.Range("M2").Value = [double]"7.341516505832146e-09"
$ws.Range("N2").Value = [double]"7.341516505832146e-09"
$ws.Range("O2").Value = 1.125815985971117
$ws.Range("P2").Value = "[0.7861843477451931, 1.46544762419704]"
$ws.Range("Q2").Value = [double]"5.634799293829929e-10"
$ws.Range("R2").Value = [double]"5.634799293829929e-10"
$ws.Range("S2").Value = 61.97368362926656
$ws.Range("T2").Value = "[54.29177437189831, 69.6555928866348]"
$ws.Range("W2").Value = 19.77357357357384
$ws.Range("X2").Value = 18.47141141141166
$ws.Range("Y2").Value = 21.07573573573602

# Row 3 updates
$ws.Range("E3").Value = 25.64000000000057
$ws.Range("G3").Value = [double]"4.440892098500626e-16"
$ws.Range("H3").Value = [double]"8.458842092382145e-15"
$ws.Range("K3").Value = 47.16565647443012
$ws.Range("L3").Value = "[33.200635281574804, 61.13067766728544]"
$ws.Range("M3").Value = [double]"1.962863205307031e-10"
$ws.Range("N3").Value = [double]"3.925726410614061e-10"
$ws.Range("O3").Value = 2.257921446724195
$ws.Range("P3").Value = "[1.9560266571900398, 2.5598162362583503]"
$ws.Range("S3").Value = 62.52674058905168
$ws.Range("T3").Value = "[55.3601833943616, 69.69329778374176]"
$ws.Range("W3").Value = 16.42602602602639
$ws.Range("X3").Value = 15.19407407407441
$ws.Range("Y3").Value = 17.65797797797838
